# Generate Report for Handback
# Rows 3 & 4 (the "a710d833..." / "d9927e3f..." files) on every sheet move
# from "Ready for handoff" to "Handed back: in sync with en-US", and the
# per-language sheets pick up their "Latest Target File" / "Latest Handback
# File" hyperlinks + a real "Latest Handback DateTime" now that the handback
# report has been generated.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# ---- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusHandedBack
$overview.Range("C3").Value = $statusHandedBack
$overview.Range("B4").Value = $statusHandedBack
$overview.Range("C4").Value = $statusHandedBack

# ---- per-language sheets --------------------------------------------
$langs = @(
    @{ Sheet = "zh-cn"; Xlf = "a710d833-2c81-4298-ad5e-1737ecc0f0ab.a1aeb35ca917befa0eb35b20a3338e348fb03ae0.zh-cn.xlf"; HandbackTime = "2016-01-25 07:59:17" },
    @{ Sheet = "de-de"; Xlf = "a710d833-2c81-4298-ad5e-1737ecc0f0ab.a1aeb35ca917befa0eb35b20a3338e348fb03ae0.de-de.xlf"; HandbackTime = "2016-01-25 07:59:38" }
)

$mdName = "a710d833-2c81-4298-ad5e-1737ecc0f0ab.md"

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)

    foreach ($r in 3, 4) {
        $ws.Range("B$r").Value = $statusHandedBack

        $ws.Hyperlinks.Add(
            $ws.Range("E$r"),
            "https://github.com/OpenLocalizationTestOrg/oltest.$($lang.Sheet)/blob/a710d833-2c81-4298-ad5e-1737ecc0f0ab/e2e/$mdName",
            [Type]::Missing,
            [Type]::Missing,
            $mdName
        ) | Out-Null

        $ws.Hyperlinks.Add(
            $ws.Range("F$r"),
            "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a710d833-2c81-4298-ad5e-1737ecc0f0ab/ol-handback/OpenLocalizationTestOrg/oltest.$($lang.Sheet)/yuwzho/$($lang.Xlf)",
            [Type]::Missing,
            [Type]::Missing,
            $lang.Xlf
        ) | Out-Null

        $ws.Range("G$r").Value = $lang.HandbackTime
    }
}

Write-Host "Done"
